$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix header label: "Maturity" -> "Time to Maturity"
$ws.Range("A1").Value = "Time to Maturity"
$ws.Range("B1").Value = "Discount Factor"

# Set column A width to fit the new, longer label (matches col width=16 bestFit)
$ws.Columns.Item(1).ColumnWidth = 15.2

# Move selection/active cell to A2
$ws.Range("A2").Select()
